$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "73.956.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +7.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.624.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +7.59%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "184.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +13.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "581.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.79%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.535"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.199"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +18.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.624.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.62%  "
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.359"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "73.796.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.36%  "
$ws.Range("E15").Value = "  +6.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.076.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +12.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.619.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +29.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +12.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.02%  "
$ws.Range("E22").Value = "  +16.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.97%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.734.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0935"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "516.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +20.28%  "
$ws.Range("E32").Value = "  +18.84%  "
$ws.Range("E33").Value = "  +6.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +12.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("E41").Value = "  +12.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.28%  "
$ws.Range("E43").Value = "  +8.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "159.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +22.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0866"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +20.45%  "
$ws.Range("E47").Value = "  +13.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "38.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.528"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +21.05%  "
